$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles + row height) of the last existing data row (38)
# down into the new row (39) before writing any values so the new row
# visually matches the rest of the table.
$ws.Range("A38:K38").Copy()
$ws.Range("A39:K39").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(39).RowHeight = $ws.Rows.Item(38).RowHeight
$excel.CutCopyMode = 0

# Fill in the new data row values
$ws.Cells.Item(39, 1).Value2 = "Ministerio de Trabajo y Desarrollo Social"
$ws.Cells.Item(39, 2).Value2 = 38
$ws.Cells.Item(39, 3).Value2 = "Trabajo"
$ws.Cells.Item(39, 4).Value2 = $ws.Cells.Item(38, 4).Value2
$ws.Cells.Item(39, 5).Value2 = "https://www.mitradel.gob.pa/mitradel-extiende-vigencia-de-permisos-de-trabajo-que-vencian-entre-marzo-y-septiembre-de-2020/"
$ws.Cells.Item(39, 6).Value2 = "El Ministerio de Trabajo y Desarrollo Laboral (Mitradel), publicó en gaceta oficial la Resolución del Decreto Ministerial 225 del 19 de agosto de 2020, por la cual se extiende la vigencia de los permisos de trabajo que vencían entre los meses de marzo a septiembre de 2020."
$ws.Cells.Item(39, 7).Value2 = "https://www.mitradel.gob.pa"
$ws.Cells.Item(39, 8).Value2 = 44064
$ws.Cells.Item(39, 9).Value2 = 44064
$ws.Cells.Item(39, 10).Value2 = "Panamá"
$ws.Cells.Item(39, 11).Value2 = "Ministerial"

# Add hyperlinks for the new row (Descarga Link column E, Sitio Web column G)
$ws.Hyperlinks.Add($ws.Cells.Item(39, 7), "https://www.mitradel.gob.pa")
$ws.Hyperlinks.Add($ws.Cells.Item(39, 5), "https://www.mitradel.gob.pa/mitradel-extiende-vigencia-de-permisos-de-trabajo-que-vencian-entre-marzo-y-septiembre-de-2020/")

# Adding a hyperlink resets the cell style to the default Hyperlink style,
# so re-apply the formatting copied from row 38 on top of those two cells.
$ws.Range("E38").Copy()
$ws.Range("E39").PasteSpecial(-4122)
$ws.Range("G38").Copy()
$ws.Range("G39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the actual values/hyperlinks text (PasteSpecial formats only touches
# formatting, so values above remain intact - re-affirm just in case)
$ws.Cells.Item(39, 5).Value2 = "https://www.mitradel.gob.pa/mitradel-extiende-vigencia-de-permisos-de-trabajo-que-vencian-entre-marzo-y-septiembre-de-2020/"
$ws.Cells.Item(39, 7).Value2 = "https://www.mitradel.gob.pa"

# Grow the data validation on column C to include the new row
$ws.Range("C2:C38").Validation.Delete()
$ws.Range("C2:C39").Validation.Add(0, 1, 1, "")
$ws.Range("C2:C39").Validation.ErrorTitle = "Entrada no válida"
$ws.Range("C2:C39").Validation.ErrorMessage = "Selecciona una categoría de la lista"
$ws.Range("C2:C39").Validation.InputTitle = "Categoria"
$ws.Range("C2:C39").Validation.InputMessage = "Selecciona una categoría de la lista"

# Grow the Excel table (ListObject) to include the new row, and resync the
# autofilter / table range accordingly.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K39"))

# Match the final selection/scroll position left by the author's edit
$null = $ws.Range("I39").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 38

Write-Output "done"
